# A1_A2.xlsx - "content+ system to progress through content"
#
# On the "A-approach" sheet (the properties table for choice A2):
#   - The "Manpower" row's value switches from the generic "null"
#     placeholder string to a concrete numeric penalty (-100).
#   - A brand new "Unique" property row is appended (row 18), mirroring
#     the formatting of the row above it ("Temp Scale"/row 17), with a
#     "null" placeholder value so it can be filled in later.
#   - The last cell clicked on the sheet (saved as the sheetView
#     selection) ends up at K32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A-approach")

# Manpower (E11) result becomes a fixed numeric value instead of "null"
$ws.Range("F11").Value = -100

# New property row: "Unique" -> "null", styled like the row above it
# (E: right aligned label, F: center aligned value)
$ws.Range("E18").Value = "Unique"
$ws.Range("E18").HorizontalAlignment = -4152  # xlRight
$ws.Range("F18").Value = "null"
$ws.Range("F18").HorizontalAlignment = -4108  # xlCenter

# Reflect the author's final selection in the saved sheet view
$ws.Range("K32").Select() | Out-Null
